$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34 (shifts existing rows 34-67 down to 35-68)
$ws.Rows("34:34").Insert()

# The new row 34 duplicates the data that was in row 33 (weekly log entry re-added)
$ws.Range("A33:R33").Copy($ws.Range("A34:R34"))
